# Corrections to the "objetos" sheet: fix imported values and wire up
# the new "lamp1" interface row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the A/C consumption value (imported incorrectly as 21).
$ws.Range("C1").Value = 16

# Fix the Televisor power value (imported incorrectly as 85).
$ws.Range("D2").Value = 0

# Connect the new "Lâmpada" (lamp1) interface as row 3.
$ws.Range("A3").Value = "lamp1"
$ws.Range("B3").Value = "Lâmpada"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = $true

# Make sure the new row uses the default (unstyled) formatting, matching
# the rest of the freshly-added data rather than inheriting the column's
# number styles.
$ws.Range("A3:D3").Style = "Normal"
